# QA_TestCase_Auto_Optimus_2_2_2 - 3 and 4
# Fill in the missing automation-run values for rows 4-7 on the
# "RelatedCounterParty" sheet, then leave that sheet active/selected
# (matching the new workbook/sheet view state), moving the active tab
# away from "Portfolio".

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("RelatedCounterParty")
$ws5 = $wb.Worksheets.Item("Portfolio")

# Row 4: add the FirstName value that was missing.
$ws3.Range("E4").Value = "Fname TXA"

# Row 5: fill in the related-party details for the delete test case.
$ws3.Range("D5").Value = "Joe Pen QA2"
$ws3.Range("E5").Value = "Fname TXA"
$ws3.Range("F5").Value = "LN Updt TXA"
$ws3.Range("G5").Value = "TXA"

# Row 6: fill in action + related-party details for the search test case.
$ws3.Range("B6").Value = "Search"
$ws3.Range("D6").Value = "Joe Pen QA2"
$ws3.Range("E6").Value = "Fname TXA"
$ws3.Range("F6").Value = "TXA"
$ws3.Range("G6").Value = "TXA"

# Row 7: fill in action + related-party details for the download test case.
$ws3.Range("B7").Value = "Download Btn"
$ws3.Range("D7").Value = "Joe Pen QA2"
$ws3.Range("E7").Value = "Fname TXA"
$ws3.Range("F7").Value = "TXA"
$ws3.Range("G7").Value = "TXA"

# Move the selection/active sheet from Portfolio to RelatedCounterParty.
$ws5.Range("E16").Select() | Out-Null

$ws3.Activate()
$ws3.Range("G8").Select() | Out-Null
